# Notenrechner Codereview SWE 2 - apply commit changes
#
# Summary of the change (see commit message / diff):
#   "Liste an FotografInnen wird angezeigt, das Fenster laesst sich
#    vergroessern und verkleinern, Details von FotografInnen wird angezeigt"
#   -> three checklist items on the "Notenrechner" sheet go from "not
#      graded" to "1/2 points awarded": B8, B9, B10, B21, B26. The B33
#      total (a SUM formula) recalculates automatically. The active
#      window/view is also resized and rescrolled/reselected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notenrechner")
$ws.Activate()

# --- Punkte (points) updates on the checklist ------------------------------
# Row 8  "IPTC & EXIF Informationen werden bei Bild angezeigt" : 0 -> 1
$ws.Cells.Item(8, 2).Value2 = 1
# Row 9  "IPTC Informationen koennen bearbeitet werden"        : 0 -> 2
$ws.Cells.Item(9, 2).Value2 = 2
# Row 10 "Liste alle Bilder ist eine UI-Komponente"             : 1 -> 2
$ws.Cells.Item(10, 2).Value2 = 2
# Row 21 "Fotografen_innen auflisten" (list photographers)      : 0 -> 1
$ws.Cells.Item(21, 2).Value2 = 1
# Row 26 "Konfiguration wird benutzt"                           : 0 -> 1
$ws.Cells.Item(26, 2).Value2 = 1

# B33 holds =SUM(B6:B27)-SUM(B30:B32); it recalculates automatically once
# the workbook is saved/recalculated, moving its cached value from 3 to 9.

# --- Window / view state -----------------------------------------------
# The author resized & repositioned the Excel window and scrolled /
# reselected the sheet (selection moves from B9 to B20, with A4 pinned as
# the top-left visible row via topLeftCell="A4"). Reproduce the same COM
# calls an interactive user would make.
$win = $wb.Windows.Item(1)
$win.Left = 1960
$win.Top = 140
$win.Width = 13720
$win.Height = 13800

$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

[void]$ws.Range("B20").Select()
